$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You are creating two tables in your database: one to store student data (e.g., class, name, and addresses) and one to store students' final exam results (e.g., class, subject, score, date_of_exam). Both the tables are connected with a common column class.Which of the following constraints should you choose for these tables?",
        "ques_type": 2,
        "options": [
            "Both primary key and foreign key",
            "Primary key only",
            "Foreign key only",
            "NOT NULL"
        ],
        "score": "Both primary key and foreign key"
    },
    {
        "title": "You are selling products that come in many colors, and you have a database containing product information, including color. You want to use a SQL query to display a list of all blue and white products. Which operator is most suitable for this query operation?",
        "ques_type": 2,
        "options": [
            "AND",
            "NOT",
            "OR",
            "&lt&gt "
        ],
        "score": "OR"
    },
    {
        "title": "You are working on a database of employee data at your company. You've been instructed to fetch the second-highest employee salary in a particular department.Which of the following must you use to write this?",
        "ques_type": 2,
        "options": [
            "Union",
            "Simple select statement",
            "Intersect",
            "Aggregate functions"
        ],
        "score": "Aggregate functions"
    },
    {
        "title": "You are working on a database containing information about students taking admissions exams. In the database, exam results update automatically in the result table whenever student scores are entered in the student table. True or false: This system uses the update trigger.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@

# Row 2 held the long questions string previously; it is being consolidated
# into A1, so clear it out entirely (collapses the used range back to A1).
$ws.Range("A2").ClearContents()

# A1 previously held a bold/bordered/centered "0" placeholder value — strip
# that formatting back to the workbook's default (Normal) style.
$ws.Range("A1").ClearFormats()

# Write the reformatted questions payload into A1.
$ws.Range("A1").Value = $questionsText
